# Helper: convert an EMU value to points for Shape.Left/Top/Width/Height,
# nudging by a hair so the engine's internal float32 storage doesn't
# truncate down to (emu-1) when it re-derives EMUs on save.
function EMUPt($emu) {
    return ($emu / 12700.0) + 0.0000125
}

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Refresh the "Last edited" date field on the slide master and
#    every slide layout (11/12/2019 -> 8/13/2025).
# ---------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        if ($shp.TextFrame.TextRange.Text -eq "11/12/2019") {
            $shp.TextFrame.TextRange.Text = "8/13/2025"
        }
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.TextRange.Text -eq "11/12/2019") {
                $shp.TextFrame.TextRange.Text = "8/13/2025"
            }
        }
    }
}

# ---------------------------------------------------------------
# 2) Revise the capability-statement diagram on slide 1:
#    merge the "CodeSystem" box into the "ValueSet" box
#    ("ValueSet, CodeSystem, etc."), drop the now-redundant
#    "CodeSystem" rectangle + its connector, and slide the
#    remaining connector down to meet the relocated box.
# ---------------------------------------------------------------
$s = $p.Slides.Item(1)

$valueSetBox = Get-ShapeByName $s "Rectangle 36"
$tr = $valueSetBox.TextFrame.TextRange
$fullText = $tr.Text
$startChar = $fullText.IndexOf("ValueSet") + 1

$boldText = "ValueSet, CodeSystem"
$etcText = ", etc."

$wordRange = $tr.Characters($startChar, "ValueSet".Length)
$wordRange.Text = $boldText + $etcText

# "ValueSet, CodeSystem" keeps the original bold run; ", etc." becomes
# a new, non-bold trailing run.
$etcRange = $tr.Characters($startChar + $boldText.Length, $etcText.Length)
$etcRange.Font.Bold = 0

# Move the merged box down into the old "CodeSystem" slot.
$valueSetBox.Left = EMUPt(962079)
$valueSetBox.Top = EMUPt(1308894)

# Remove the "CodeSystem" rectangle and its connector arrow.
$codeSystemBox = Get-ShapeByName $s "Rectangle 7"
if ($codeSystemBox -ne $null) { $codeSystemBox.Delete() }

$lowerConnector = Get-ShapeByName $s "Straight Arrow Connector 8"
if ($lowerConnector -ne $null) { $lowerConnector.Delete() }

# Reposition the remaining connector so it still points at the box.
$upperConnector = Get-ShapeByName $s "Straight Arrow Connector 37"
$upperConnector.Left = EMUPt(952592)
$upperConnector.Top = EMUPt(1156494)
